$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Modelo: AREPD)
$ws.Cells.Item(2, 7).Value  = 97.75158429315461
$ws.Cells.Item(2, 8).Value  = 97.24752963974205
$ws.Cells.Item(2, 9).Value  = 95.90589757124256

# Row 3 (Modelo: AV-MCPS)
$ws.Cells.Item(3, 2).Value  = 76.76783679102182
$ws.Cells.Item(3, 3).Value  = 94.09497597764282
$ws.Cells.Item(3, 4).Value  = 98.9975916113487
$ws.Cells.Item(3, 5).Value  = 98.71559333961997
$ws.Cells.Item(3, 6).Value  = 98.45944038656165
$ws.Cells.Item(3, 7).Value  = 97.99313915062997
$ws.Cells.Item(3, 8).Value  = 97.51695878566298
$ws.Cells.Item(3, 9).Value  = 96.17801537849368

# Row 4 (Modelo: Block Bootstrapping)
$ws.Cells.Item(4, 2).Value  = 91.92663353917045
$ws.Cells.Item(4, 3).Value  = 98.13286716457354
$ws.Cells.Item(4, 4).Value  = 99.23769815596334
$ws.Cells.Item(4, 5).Value  = 98.70980009956406
$ws.Cells.Item(4, 6).Value  = 98.21291579568853
$ws.Cells.Item(4, 7).Value  = 97.73742794960029
$ws.Cells.Item(4, 8).Value  = 97.23354418526307
$ws.Cells.Item(4, 9).Value  = 95.89844826086346

# Row 5 (Modelo: DeepAR)
$ws.Cells.Item(5, 7).Value  = 97.82267973387349
$ws.Cells.Item(5, 8).Value  = 97.29391776958097
$ws.Cells.Item(5, 9).Value  = 95.98217387238336

# Row 6 (Modelo: EnCQR-LSTM)
$ws.Cells.Item(6, 2).Value  = 86.45074642503886
$ws.Cells.Item(6, 3).Value  = 96.55737218175472
$ws.Cells.Item(6, 4).Value  = 98.9514846297614
$ws.Cells.Item(6, 5).Value  = 98.61262141884974
$ws.Cells.Item(6, 6).Value  = 98.4147800531138
$ws.Cells.Item(6, 7).Value  = 97.84263269141981
$ws.Cells.Item(6, 8).Value  = 97.4998953706876
$ws.Cells.Item(6, 9).Value  = 96.09570157831996

# Row 7 (Modelo: LSPM)
$ws.Cells.Item(7, 7).Value  = 98.03927506433419
$ws.Cells.Item(7, 8).Value  = 97.60496829686596
$ws.Cells.Item(7, 9).Value  = 96.28492872810629

# Row 8 (Modelo: LSPMW)
$ws.Cells.Item(8, 7).Value  = 98.0263346426579
$ws.Cells.Item(8, 8).Value  = 97.53753934234044
$ws.Cells.Item(8, 9).Value  = 96.21960877054136

# Row 9 (Modelo: MCPS)
$ws.Cells.Item(9, 2).Value  = 74.50297034586364
$ws.Cells.Item(9, 3).Value  = 95.16865857826339
$ws.Cells.Item(9, 4).Value  = 99.37011432484057
$ws.Cells.Item(9, 5).Value  = 98.88479443372457
$ws.Cells.Item(9, 6).Value  = 98.40383018180337
$ws.Cells.Item(9, 7).Value  = 97.94108433439878
$ws.Cells.Item(9, 8).Value  = 97.43003042753313
$ws.Cells.Item(9, 9).Value  = 96.0692887529015

# Row 10 (Modelo: Sieve Bootstrap)
$ws.Cells.Item(10, 7).Value = 98.09893715101707
$ws.Cells.Item(10, 8).Value = 98.38694656126903
$ws.Cells.Item(10, 9).Value = 95.47756413258084
